# Update excel models and documentation
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Shared-string text edit ---------------------------------------------
# D9 ("TaskName" for row 9) held the string "Analysis"; rename it in place.
$ws.Range("D9").Value = "Backend Analysis"

# --- New TaskDependencies (column F) values on rows 3, 4 and 6 ----------
$ws.Range("F3").Value = 1
$ws.Range("F4").Value = 2
$ws.Range("F6").Value = 4

# --- Rows 5 & 6: TaskName / EstimatedEffortHours / InternalID swap ------
$ws.Range("D5").Value = "Database Setup"
$ws.Range("E5").Value = 60
$ws.Range("I5").Value = 1327

$ws.Range("D6").Value = "API Development"
$ws.Range("E6").Value = 50
$ws.Range("I6").Value = 1236

# --- EstimatedEffortHours swap between rows 9 and 11 ---------------------
$ws.Range("E9").Value = 47
$ws.Range("E11").Value = 45

# --- View state: zoom + active selection ---------------------------------
$excel.ActiveWindow.Zoom = 70
[void]$ws.Range("G7").Select()
